# Generate Report for Handoff
# Updates the localization-status workbook to reflect that b.md has been
# re-handed-off (new handoff xliff files were generated), while a.md is
# unchanged / remains "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$newStatus       = "Ready for handoff"
$newHoDate       = "2016-08-30 14:48:43"
$zhHandoffFile   = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate   = "2016-08-30 14:48:39"
$deHandoffFile   = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate   = "2016-08-30 14:48:43"
$errorDetail     = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/61502a264614e7763592cf91f36bf31cf6c824b1/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/525b2bbdcdbb784024a74fca961f5dd17c5b3f82/e2e/b.md."

# --- Overview sheet: row 3 is the b.md entry ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Range("G3").Value = $newHoDate

# --- zh-cn sheet: row 3 is the b.md entry ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
# Copy the (text) "False" from O3 so the cell keeps its text type instead of
# being auto-coerced into a boolean by a plain string assignment.
$zhcn.Range("O3").Copy($zhcn.Range("F3"))
$zhcn.Range("G3").Value = $zhHandoffFile
$zhcn.Range("H3").Value = $zhHandoffDate
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Range("P:P").ColumnWidth = 39.14

# --- de-de sheet: row 3 is the b.md entry ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("O3").Copy($dede.Range("F3"))
$dede.Range("G3").Value = $deHandoffFile
$dede.Range("H3").Value = $deHandoffDate
$dede.Range("P3").Value = $errorDetail
$dede.Range("P:P").ColumnWidth = 39.14
